$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure new time cells (column B) inherit the existing time number format
# used by B2 (h:mm:ss, style index 2 in the original workbook).
$ws.Range("B3:B9").NumberFormat = $ws.Range("B2").NumberFormat

# Row 3
$ws.Range("A3").Value = "30-09-2025"
$ws.Range("B3").Value = 0.64778935185185182
$ws.Range("C3").Value = "Rangeen Pop"
$ws.Range("D3").Value = "Outer"
$ws.Range("E3").Value = 10
$ws.Range("F3").Value = 500
$ws.Range("G3").Value = 456
$ws.Range("H3").Value = 2
$ws.Range("I3").Value = "Credit"
$ws.Range("J3").Value = 50
$ws.Range("K3").Value = 12

# Row 4
$ws.Range("A4").Value = "30-09-2025"
$ws.Range("B4").Value = 0.6731018518518519
$ws.Range("C4").Value = "Rangeen Pop"
$ws.Range("D4").Value = "Outer"
$ws.Range("E4").Value = 10
$ws.Range("F4").Value = 500
$ws.Range("G4").Value = 456
$ws.Range("H4").Value = 5
$ws.Range("I4").Value = "Online"
$ws.Range("J4").Value = 50

# Row 5
$ws.Range("A5").Value = "30-09-2025"
$ws.Range("B5").Value = 0.68487268518518518
$ws.Range("C5").Value = "Rangeen Pop"
$ws.Range("D5").Value = "Outer"
$ws.Range("E5").Value = 10
$ws.Range("F5").Value = 500
$ws.Range("G5").Value = 75
$ws.Range("H5").Value = 2
$ws.Range("I5").Value = "Online"
$ws.Range("J5").Value = 50
$ws.Range("K5").Value = 50

# Row 6
$ws.Range("A6").Value = "30-09-2025"
$ws.Range("B6").Value = 0.76513888888888892
$ws.Range("C6").Value = "Rangeen Pop"
$ws.Range("D6").Value = "Outer"
$ws.Range("E6").Value = 10
$ws.Range("F6").Value = 50
$ws.Range("G6").Value = 1000
$ws.Range("H6").Value = 456
$ws.Range("I6").Value = 5
$ws.Range("J6").Value = "Online"

# Row 7
$ws.Range("A7").Value = "30-09-2025"
$ws.Range("B7").Value = 0.76613425925925926
$ws.Range("C7").Value = "Jalebi Jelly"
$ws.Range("D7").Value = "Outer"
$ws.Range("E7").Value = 5
$ws.Range("F7").Value = 60
$ws.Range("G7").Value = 900
$ws.Range("H7").Value = 654
$ws.Range("I7").Value = 2
$ws.Range("J7").Value = "Credit"

# Row 8
$ws.Range("A8").Value = "30-09-2025"
$ws.Range("B8").Value = 0.77379629629629632
$ws.Range("C8").Value = "Jalebi Jelly"
$ws.Range("D8").Value = "Outer"
$ws.Range("E8").Value = 20
$ws.Range("F8").Value = 50
$ws.Range("G8").Value = 1250
$ws.Range("H8").Value = 46
$ws.Range("I8").Value = 3
$ws.Range("J8").Value = "Online"

# Row 9
$ws.Range("A9").Value = "30-09-2025"
$ws.Range("B9").Value = 0.77614583333333331
$ws.Range("C9").Value = "Jalebi Jelly"
$ws.Range("D9").Value = "Outer"
$ws.Range("E9").Value = 5
$ws.Range("F9").Value = 50
$ws.Range("G9").Value = 1250
$ws.Range("H9").Value = "sc"
$ws.Range("I9").Value = 3
$ws.Range("J9").Value = "Credit"
